$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: rename columns to match new plan-import schema
$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

# Move the active cell selection to F2
$ws.Range("F2").Select()
